$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural: make room for the new "Power 90" block at the top,
#     and for the extra blank separator row above the "Power 50" block ---
$ws.Rows("1:3").Insert()
$ws.Rows("12:12").Insert()

# --- Row 1/2: new "Power 90" section ---
$ws.Range("A1").Value = "Power 90"
$ws.Range("B1").Value = -180
$ws.Range("C1").Value = -90
$ws.Range("D1").Value = -45
$ws.Range("E1").Value = -30
$ws.Range("H1").Value = 180
$ws.Range("I1").Value = 90
$ws.Range("J1").Value = 45
$ws.Range("K1").Value = 30

$ws.Range("B2").Value = -177
$ws.Range("H2").Value = 181.2
$ws.Range("I2").Value = 88.3
$ws.Range("J2").Value = 43.2
$ws.Range("K2").Value = 26.6

# --- Row 4/5: "Power 80" section (was row 1/2) ---
$ws.Range("A4").Value = "Power 80"
$ws.Range("F4").ClearContents()

$ws.Range("B5").Value = -179
$ws.Range("C5").Value = -86.4
$ws.Range("D5").Value = -41
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()

# --- Row 7/8: "Power 70" section (was row 4/5), now with H:K columns ---
$ws.Range("A7").Value = "Power 70"
$ws.Range("B7").Value = -180
$ws.Range("C7").Value = -90
$ws.Range("D7").Value = -45
$ws.Range("E7").Value = -30
$ws.Range("F7").ClearContents()
$ws.Range("H7").Value = 180
$ws.Range("I7").Value = 90
$ws.Range("J7").Value = 45
$ws.Range("K7").Value = 30

$ws.Range("B8").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").Value = -26.2
$ws.Range("F8").ClearContents()
$ws.Range("H8").Value = 181.2
$ws.Range("I8").Value = 88.4
$ws.Range("J8").Value = 43
$ws.Range("K8").Value = 27

# --- Row 10: "Power 60" section (was row 7) ---
$ws.Range("A10").Value = "Power 60"

# --- Rows 13/14: "Power 50" section (was row 9/10) - values unchanged ---
$ws.Range("A13").Value = "Power 50"

# --- Selection matches the authored state ---
$ws.Range("C5").Select()
